# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) on each job sheet with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3371.7856
$ws.Cells.Item(76, 9).Value = 3489.4443
$ws.Cells.Item(76, 11).Value = 3489.4443
$ws.Cells.Item(76, 13).Value = -3174.4443
$ws.Cells.Item(79, 8).Value = 3371.7856
$ws.Cells.Item(79, 9).Value = 3489.4443
$ws.Cells.Item(79, 11).Value = 3489.4443
$ws.Cells.Item(79, 13).Value = -2397.4443
$ws.Cells.Item(137, 8).Value = 2502196.8
$ws.Cells.Item(137, 9).Value = 2942808.5
$ws.Cells.Item(137, 11).Value = 8828425.5
$ws.Cells.Item(137, 13).Value = -8825875.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1391.2424
$ws.Cells.Item(45, 9).Value = 1289.2963
$ws.Cells.Item(45, 11).Value = 1289.2963
$ws.Cells.Item(45, 13).Value = -912.2963
$ws.Cells.Item(74, 8).Value = 4501653.5
$ws.Cells.Item(74, 9).Value = 5578860.5
$ws.Cells.Item(74, 11).Value = 5578860.5
$ws.Cells.Item(74, 13).Value = -5577986.5
$ws.Cells.Item(77, 8).Value = 4501653.5
$ws.Cells.Item(77, 9).Value = 5578860.5
$ws.Cells.Item(77, 11).Value = 27894302.5
$ws.Cells.Item(77, 13).Value = -27889934.5
$ws.Cells.Item(102, 8).Value = 14287905
$ws.Cells.Item(102, 9).Value = 23811340
$ws.Cells.Item(102, 11).Value = 23811340
$ws.Cells.Item(102, 13).Value = -23809718
$ws.Cells.Item(132, 8).Value = 73526.586
$ws.Cells.Item(132, 9).Value = 46330.59
$ws.Cells.Item(132, 11).Value = 138991.77
$ws.Cells.Item(132, 13).Value = -136461.77
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(134, 8).Value = 2675.328
$ws.Cells.Item(134, 9).Value = 2350.818
$ws.Cells.Item(134, 10).Value = 5650
$ws.Cells.Item(134, 11).Value = 7052.454000000001
$ws.Cells.Item(134, 12).Value = 16950
$ws.Cells.Item(134, 13).Value = -4517.454000000001
$ws.Cells.Item(134, 14).Value = -22020
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2647.6042
$ws.Cells.Item(31, 9).Value = 1773.8077
$ws.Cells.Item(31, 10).Value = 3680.2727
$ws.Cells.Item(31, 11).Value = 1773.8077
$ws.Cells.Item(31, 12).Value = 3680.2727
$ws.Cells.Item(31, 13).Value = -1478.8077
$ws.Cells.Item(31, 14).Value = -4270.2727
$ws.Cells.Item(34, 8).Value = 2647.6042
$ws.Cells.Item(34, 9).Value = 1773.8077
$ws.Cells.Item(34, 10).Value = 3680.2727
$ws.Cells.Item(34, 11).Value = 1773.8077
$ws.Cells.Item(34, 12).Value = 3680.2727
$ws.Cells.Item(34, 13).Value = -1571.8077
$ws.Cells.Item(34, 14).Value = -4084.2727
$ws.Cells.Item(125, 8).Value = 35000
$ws.Cells.Item(125, 10).Value = 35000
$ws.Cells.Item(125, 12).Value = 35000
$ws.Cells.Item(125, 14).Value = -39920
$ws.Cells.Item(132, 8).Value = 23807
$ws.Cells.Item(132, 9).Value = 1393.5
$ws.Cells.Item(132, 10).Value = 145480.28
$ws.Cells.Item(132, 11).Value = 4180.5
$ws.Cells.Item(132, 12).Value = 436440.84
$ws.Cells.Item(132, 13).Value = -1650.5
$ws.Cells.Item(132, 14).Value = -441500.84
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(94, 8).Value = 2770
$ws.Cells.Item(126, 8).Value = 1916.6666
$ws.Cells.Item(126, 10).Value = 2180
$ws.Cells.Item(126, 12).Value = 6540
$ws.Cells.Item(126, 14).Value = -16420
$ws.Cells.Item(131, 8).Value = 1306.8182
$ws.Cells.Item(131, 10).Value = 1407.8948
$ws.Cells.Item(131, 12).Value = 4223.6844
$ws.Cells.Item(131, 14).Value = -14303.6844
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3772.9167
$ws.Cells.Item(80, 9).Value = 2880
$ws.Cells.Item(80, 10).Value = 4070.5557
$ws.Cells.Item(80, 11).Value = 2880
$ws.Cells.Item(80, 12).Value = 4070.5557
$ws.Cells.Item(80, 13).Value = -1882
$ws.Cells.Item(80, 14).Value = -6066.5557
$ws.Cells.Item(83, 8).Value = 3772.9167
$ws.Cells.Item(83, 9).Value = 2880
$ws.Cells.Item(83, 10).Value = 4070.5557
$ws.Cells.Item(83, 11).Value = 14400
$ws.Cells.Item(83, 12).Value = 20352.7785
$ws.Cells.Item(83, 13).Value = -9408
$ws.Cells.Item(83, 14).Value = -30336.7785
$ws.Cells.Item(113, 8).Value = 1782.2
$ws.Cells.Item(113, 9).Value = 1677.75
$ws.Cells.Item(113, 10).Value = 2200
$ws.Cells.Item(113, 11).Value = 1677.75
$ws.Cells.Item(113, 12).Value = 2200
$ws.Cells.Item(113, 13).Value = 492.25
$ws.Cells.Item(113, 14).Value = -6540
$ws.Cells.Item(126, 8).Value = 2204
$ws.Cells.Item(126, 9).Value = 2497.4
$ws.Cells.Item(126, 10).Value = 1959.5
$ws.Cells.Item(126, 11).Value = 7492.200000000001
$ws.Cells.Item(126, 12).Value = 5878.5
$ws.Cells.Item(126, 13).Value = -5022.200000000001
$ws.Cells.Item(126, 14).Value = -10818.5
$ws.Cells.Item(132, 8).Value = 88686.17
$ws.Cells.Item(132, 9).Value = 51489.2
$ws.Cells.Item(132, 11).Value = 154467.6
$ws.Cells.Item(132, 13).Value = -151937.6
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2879.6
$ws.Cells.Item(7, 9).Value = 2700
$ws.Cells.Item(7, 10).Value = 2924.5
$ws.Cells.Item(7, 11).Value = 2700
$ws.Cells.Item(7, 12).Value = 2924.5
$ws.Cells.Item(7, 13).Value = -2588
$ws.Cells.Item(7, 14).Value = -3148.5
$ws.Cells.Item(61, 8).Value = 2119
$ws.Cells.Item(61, 9).Value = 1708.6666
$ws.Cells.Item(61, 10).Value = 3350
$ws.Cells.Item(61, 11).Value = 1708.6666
$ws.Cells.Item(61, 12).Value = 3350
$ws.Cells.Item(61, 13).Value = -1506.6666
$ws.Cells.Item(61, 14).Value = -3754
$ws.Cells.Item(113, 8).Value = 2119
$ws.Cells.Item(113, 9).Value = 1708.6666
$ws.Cells.Item(113, 10).Value = 3350
$ws.Cells.Item(113, 11).Value = 1708.6666
$ws.Cells.Item(113, 12).Value = 3350
$ws.Cells.Item(113, 13).Value = 461.3334
$ws.Cells.Item(113, 14).Value = -7690
$ws.Cells.Item(126, 8).Value = 2879.6
$ws.Cells.Item(126, 9).Value = 2700
$ws.Cells.Item(126, 10).Value = 2924.5
$ws.Cells.Item(126, 11).Value = 8100
$ws.Cells.Item(126, 12).Value = 8773.5
$ws.Cells.Item(126, 13).Value = -5630
$ws.Cells.Item(126, 14).Value = -13713.5
$ws.Cells.Item(127, 8).Value = 49966.668
$ws.Cells.Item(127, 10).Value = 49966.668
$ws.Cells.Item(127, 12).Value = 49966.668
$ws.Cells.Item(127, 14).Value = -59886.668
$ws.Cells.Item(136, 8).Value = 61618.09
$ws.Cells.Item(136, 9).Value = 34700.5
$ws.Cells.Item(136, 11).Value = 104101.5
$ws.Cells.Item(136, 13).Value = -101551.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 53301.9
$ws.Cells.Item(136, 9).Value = 38547.184
$ws.Cells.Item(136, 10).Value = 86500
$ws.Cells.Item(136, 11).Value = 115641.552
$ws.Cells.Item(136, 12).Value = 259500
$ws.Cells.Item(136, 13).Value = -113091.552
$ws.Cells.Item(136, 14).Value = -264600
$ws.Cells.Item(137, 8).Value = 67500
$ws.Cells.Item(137, 10).Value = 67500
$ws.Cells.Item(137, 12).Value = 67500
$ws.Cells.Item(137, 14).Value = -77700
